$d = $word.ActiveDocument

# 1. Fix the "spealer" typo -> "speaker" in the PRESCRIPTIVIST answer cell.
$d.Content.Find.Execute("spealer", $true, $false, $false, $false, $false,
                         $true, 1, $false, "speaker", 2) | Out-Null

# 2. Normalise the answer-key tables' indentation / padding / column width
#    (tblInd 55->53 dxa, cell left-margin 53->50 dxa, 2nd column 9435->9434 dxa)
#    across every table in the document.
for ($ti = 1; $ti -le $d.Tables.Count; $ti++) {
    $t = $d.Tables.Item($ti)

    $t.Rows.LeftIndent = 2.65       # 53 dxa
    $t.LeftPadding = 2.5            # 50 dxa
    if ($t.Columns.Count -ge 2) {
        $t.Columns.Item(2).Width = 471.7   # 9434 dxa
    }

    for ($r = 1; $r -le $t.Rows.Count; $r++) {
        for ($c = 1; $c -le $t.Columns.Count; $c++) {
            try {
                $cell = $t.Cell($r, $c)
                $cell.LeftPadding = 2.5
            } catch {
                # vertically-merged continuation cells aren't reachable via
                # Table.Cell(); fall back to the row's cell collection.
                $row = $t.Rows.Item($r)
                foreach ($rc in $row.Cells) {
                    if ($rc.ColumnIndex -eq $c) {
                        $rc.LeftPadding = 2.5
                    }
                }
            }
        }
    }
}

# 3. Re-apply the Normal style's paragraph/run formatting so it matches the
#    canonical re-save (explicit left-to-right / left-aligned paragraph
#    formatting and the resolved "automatic" text colour).
$normal = $d.Styles.Item("Normal")
$normal.ParagraphFormat.ReadingOrder = 0
$normal.ParagraphFormat.Alignment = 0
$normal.Font.Color = 655360
